# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") held stale "Strike#" derived values; recompute them
# and rewrite the column for every data row (rows 2-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, keyed by row number, as produced by the regenerated save_data.
$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 2
    6  = 1
    7  = 2
    8  = 3
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 3
    14 = 1
    15 = 3
    16 = 4
    17 = 3
    18 = 5
    19 = 1
    20 = 8
    21 = 3
    22 = 5
    23 = 4
    24 = 6
    25 = 4
    26 = 4
    27 = 1
    28 = 3
    29 = 4
    30 = 3
    31 = 3
    32 = 1
    33 = 3
    34 = 4
    35 = 5
    36 = 5
    37 = 4
    38 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
